$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A142").Value = "IMX-USD"
$ws.Range("A143").Value = "GRT-USD"
